# Update the "Förändrad" (Changed) column C for all data rows (2-130)
# from serial date 45190 (2023-09-21) to serial date 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 130 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
